$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Disease Ontology (DO) version in E3
$ws.Range("E3").Value = "v2025-08-01"

# Update Experimental Factor Ontology (EFO) version in E4
$ws.Range("E4").Value = "v3.80.0"

# Update the active selection to match the saved cursor position
$ws.Range("E4").Select()
